$d = $word.ActiveDocument
$endPos = $d.Content.End
$insPoint = $d.Range($endPos, $endPos)
$xmlFrag = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pBdr/><w:contextualSpacing w:val="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p><w:p><w:pPr><w:pBdr/><w:contextualSpacing w:val="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p><w:p><w:pPr><w:pBdr/><w:contextualSpacing w:val="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading2"/><w:pBdr/><w:contextualSpacing w:val="0"/><w:rPr/></w:pPr><w:bookmarkStart w:colFirst="0" w:colLast="0" w:name="_qsw5bxbd7gid" w:id="3"/><w:bookmarkEnd w:id="3"/><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Considerations for Discussion</w:t></w:r></w:p><w:p><w:pPr><w:pBdr/><w:contextualSpacing w:val="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:cs="Arial Unicode MS" w:eastAsia="Arial Unicode MS" w:hAnsi="Arial Unicode MS"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">How do the associations I found for SES -&gt; sleep, SES -&gt; CRP, and sleep → CRP compare to the literature? If I didn’t find as big of an effect, would a bigger effect in other studies indicate higher likelihood of mediation?</w:t></w:r></w:p><w:p><w:pPr><w:pBdr/><w:contextualSpacing w:val="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p><w:p><w:pPr><w:pBdr/><w:contextualSpacing w:val="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">In studies finding mediation by BMI and physical activity and smoking, was the association fully explained or was there an unexplained portion? If so, hypothesize other mediators that might be present?</w:t></w:r></w:p><w:p><w:pPr><w:pBdr/><w:contextualSpacing w:val="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p><w:p><w:pPr><w:pBdr/><w:contextualSpacing w:val="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Effect of participation rate, selection bias? Was there social desirability bias in sleep questions? What about income and education questions?</w:t></w:r></w:p><w:p><w:pPr><w:pBdr/><w:contextualSpacing w:val="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p><w:p><w:pPr><w:pBdr/><w:contextualSpacing w:val="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Very little missing data for main variables, only for confounders. However, confounders may be poorly controlled as a result. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$null = $insPoint.InsertXML($xmlFrag)
Write-Host ("New paragraph count: " + $d.Paragraphs.Count)
